$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "needed" placeholder values in column F (rows 2-4) with "placeholder"
$ws.Range("F2").Value = "placeholder"
$ws.Range("F3").Value = "placeholder"
$ws.Range("F4").Value = "placeholder"

# Match the last-selected cell recorded in the saved view state
$ws.Range("F4").Select()
